$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.211.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.47%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.830.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.70%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'0.9998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.11%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'237.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.07%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'0.6066"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -3.66%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  +0.11%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.07094"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -4.68%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = "'  -2.67%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'24.03"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -3.26%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.07642"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.24%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'1.830.43"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.93%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'4.803"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -3.60%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.6380"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -5.96%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.000009964"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -2.65%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'2.067.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.65%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'79.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.70%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'5.993"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -4.27%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'29.175.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.51%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'230.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.56%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "'11.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -4.07%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.14%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'7.011"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -5.23%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'1.001"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.08%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'155.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.66%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'8.063"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -4.95%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'0.1293"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -4.38%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'16.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -3.86%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'0.06637"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.50%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'1.450"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.38%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'1.458"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.91%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'3.835"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -5.41%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'3.816"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -6.13%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'1.130"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.77%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'1.722"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -6.32%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'0.6577"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -5.57%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D38").Value = "'1.235.37"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.69%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'2.756"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.03%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.01767"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -4.61%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'6.592"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -2.99%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.9277"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.27%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").Value = "'  +0.14%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'1.982.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.16%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'100.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.50%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'63.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -3.10%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'  -2.03%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'1.634"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -4.46%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'8.512"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -5.60%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = "'  -1.47%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.1083"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -5.40%  "
$ws.Range("E51").Style = "Normal"
